# Adds the bold, red "05.04" date to the two empty "first session date"
# cells belonging to the Л17 and Л18 lecture rows of the schedule table
# (table 1). The Л17 row's empty paragraph additionally gains <w:b/> and
# <w:color w:val="FF0000"/> on its own (paragraph-mark) rPr, matching the
# committed diff.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Add-RedBoldDate($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)

    # Give the (currently empty) paragraph mark the bold + red formatting
    # before inserting text, so <w:pPr><w:rPr> picks it up too.
    $markRange = $cell.Range.Paragraphs.Item(1).Range
    $markRange.Font.Bold = $true
    $markRange.Font.Color = 255

    # Insert the date text into the (still empty) cell.
    $insertRange = $cell.Range
    $insertRange.InsertBefore($text)

    # Re-fetch the cell/range - the previous Range handles are stale once
    # the content changes - then format just the inserted characters.
    $cell2 = $table.Cell($row, $col)
    $start = $cell2.Range.Start
    $textRange = $d.Range($start, $start + $text.Length)
    $textRange.Font.Name = "Times New Roman"
    $textRange.Font.Size = 14
    $textRange.Font.Bold = $true
    $textRange.Font.Color = 255
}

# Row 32 = "Л17" lecture row, column 2 = first ("занять") date column.
Add-RedBoldDate $t 32 2 "05.04"

# Row 33 = "Л18" lecture row, column 2 = first ("занять") date column.
Add-RedBoldDate $t 33 2 "05.04"
